$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.722.15'
$ws.Range('E2').Value = '  +2.29%  '
$ws.Range('D3').Value = '3.550.17'
$ws.Range('E3').Value = '  +0.77%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '581.17'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.44%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '185.01'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.62%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.629'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.99%  '
$ws.Range('D8').Value = '3.537.65'
$ws.Range('E8').Value = '  +0.55%  '
$ws.Range('E9').Value = '  +0.03%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.220'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +18.64%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.651'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.53%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '54.28'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.58%  '
$ws.Range('E13').Value = '  +4.12%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.47'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.11%  '
$ws.Range('D15').Value = '4.118.12'
$ws.Range('E15').Value = '  +0.73%  '
$ws.Range('D16').Value = '70.780.96'
$ws.Range('E16').Value = '  +2.41%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '19.20'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.86%  '
$ws.Range('D18').Value = '3.547.20'
$ws.Range('E18').Value = '  +0.88%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '570.94'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +5.97%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.35'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.64%  '
$ws.Range('E21').Value = '  +0.42%  '
$ws.Range('E22').Value = '  -3.01%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '17.67'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -14.54%  '
$ws.Range('B24').Value = 'PancakeSwap'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.54'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.64%  '
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.97'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.04%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '95.52'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.79%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.24'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.36%  '
$ws.Range('E28').Value = '  +0.34%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.13'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.56%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '32.41'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.57%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.19'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.88%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '12.25'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.48%  '
$ws.Range('E33').Value = '  +2.27%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '63.25'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.28%  '
$ws.Range('E35').Value = '  +11.56%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '544.11'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -5.09%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.413'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.01%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.38'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +8.45%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '37.72'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.99%  '
$ws.Range('E40').Value = '  +0.02%  '
$ws.Range('D41').Value = '0.0₃0797'
$ws.Range('E41').Value = '  +4.04%  '
$ws.Range('D42').Value = '3.573.38'
$ws.Range('E42').Value = '  +11.39%  '
$ws.Range('E43').Value = '  +1.62%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.41'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.89%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0452'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.53%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.49'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.79%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.91'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.07%  '
$ws.Range('E48').Value = '  +0.93%  '
$ws.Range('E49').Value = '  +2.28%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.50'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +12.33%  '
$ws.Range('B51').Value = 'FLOKI'
$ws.Range('C51').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.000263'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +15.30%  '
